{"js": "const replacements = [\n  [\"2023-07-21 Friday\", \"2023-07-22 Saturday\"],\n  [\"30\u00f78=\", \"29\u00f74=\"],\n  [\"11\u00f77=\", \"10\u00f77=\"],\n  [\"80\u00f73=\", \"55\u00f77=\"],\n  [\"45\u00f73=\", \"40\u00f78=\"],\n  [\"46\u00f75=\", \"42\u00f77=\"],\n  [\"18\u00f77=\", \"92\u00f74=\"],\n  [\"56\u00f73=\", \"49\u00f72=\"],\n  [\"20\u00f72=\", \"51\u00f76=\"],\n  [\"83\u00f75=\", \"24\u00f79=\"],\n  [\"50\u00f72=\", \"22\u00f72=\"],\n  [\"92\u00f78=\", \"64\u00f73=\"],\n  [\"88\u00f76=\", \"82\u00f78=\"],\n  [\"91\u00f76=\", \"51\u00f77=\"],\n  [\"19\u00f77=\", \"54\u00f78=\"],\n  [\"43\u00f77=\", \"90\u00f75=\"],\n  [\"32\u00f79=\", \"66\u00f74=\"],\n  [\"78\u00f77=\", \"25\u00f75=\"],\n  [\"98\u00f72=\", \"36\u00f78=\"],\n  [\"24\u00f74=\", \"49\u00f73=\"],\n  [\"89\u00f78=\", \"78\u00f75=\"],\n  [\"37\u00f78=\", \"97\u00f77=\"],\n  [\"38\u00f74=\", \"35\u00f75=\"],\n  [\"42\u00f76=\", \"65\u00f76=\"],\n  [\"67\u00f76=\", \"19\u00f76=\"],\n  [\"72\u00f76=\", \"85\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every division problem in the practice table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2023-07-21 Friday\"; New = \"2023-07-22 Saturday\" },\n    @{ Old = \"30\u00f78=\";             New = \"29\u00f74=\" },\n    @{ Old = \"11\u00f77=\";             New = \"10\u00f77=\" },\n    @{ Old = \"80\u00f73=\";             New = \"55\u00f77=\" },\n    @{ Old = \"45\u00f73=\";             New = \"40\u00f78=\" },\n    @{ Old = \"46\u00f75=\";             New = \"42\u00f77=\" },\n    @{ Old = \"18\u00f77=\";             New = \"92\u00f74=\" },\n    @{ Old = \"56\u00f73=\";             New = \"49\u00f72=\" },\n    @{ Old = \"20\u00f72=\";             New = \"51\u00f76=\" },\n    @{ Old = \"83\u00f75=\";             New = \"24\u00f79=\" },\n    @{ Old = \"50\u00f72=\";             New = \"22\u00f72=\" },\n    @{ Old = \"92\u00f78=\";             New = \"64\u00f73=\" },\n    @{ Old = \"88\u00f76=\";             New = \"82\u00f78=\" },\n    @{ Old = \"91\u00f76=\";             New = \"51\u00f77=\" },\n    @{ Old = \"19\u00f77=\";             New = \"54\u00f78=\" },\n    @{ Old = \"43\u00f77=\";             New = \"90\u00f75=\" },\n    @{ Old = \"32\u00f79=\";             New = \"66\u00f74=\" },\n    @{ Old = \"78\u00f77=\";             New = \"25\u00f75=\" },\n    @{ Old = \"98\u00f72=\";             New = \"36\u00f78=\" },\n    @{ Old = \"24\u00f74=\";             New = \"49\u00f73=\" },\n    @{ Old = \"89\u00f78=\";             New = \"78\u00f75=\" },\n    @{ Old = \"37\u00f78=\";             New = \"97\u00f77=\" },\n    @{ Old = \"38\u00f74=\";             New = \"35\u00f75=\" },\n    @{ Old = \"42\u00f76=\";             New = \"65\u00f76=\" },\n    @{ Old = \"67\u00f76=\";             New = \"19\u00f76=\" },\n    @{ Old = \"72\u00f76=\";             New = \"85\u00f74=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Execute($r.Old, $false, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
